$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2447.9395
$ws.Range("I15").Value = 2447.9395
$ws.Range("K15").Value = 7343.818499999999
$ws.Range("M15").Value = -7174.818499999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 716043.5600000001
$ws.Range("I19").Value = 1251454.8
$ws.Range("J19").Value = 2162
$ws.Range("K19").Value = 1251454.8
$ws.Range("L19").Value = 2162
$ws.Range("M19").Value = -1251279.8
$ws.Range("N19").Value = -2512

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 920.26666
$ws.Range("I41").Value = 1055.5454
$ws.Range("J41").Value = 548.25
$ws.Range("K41").Value = 1055.5454
$ws.Range("L41").Value = 548.25
$ws.Range("M41").Value = -615.5454
$ws.Range("N41").Value = -1428.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4722.9414
$ws.Range("J51").Value = 4229.6
$ws.Range("L51").Value = 4229.6
$ws.Range("N51").Value = -5197.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 14201.8
$ws.Range("I70").Value = 7040
$ws.Range("K70").Value = 21120
$ws.Range("M70").Value = -20850

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 14201.8
$ws.Range("I73").Value = 7040
$ws.Range("K73").Value = 21120
$ws.Range("M73").Value = -20184

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 838.4375
$ws.Range("I86").Value = 843.36365
$ws.Range("J86").Value = 827.6
$ws.Range("K86").Value = 843.36365
$ws.Range("L86").Value = 827.6
$ws.Range("M86").Value = 279.63635
$ws.Range("N86").Value = -3073.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 28450
$ws.Range("J87").Value = 28450
$ws.Range("L87").Value = 28450
$ws.Range("N87").Value = -30946

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 838.4375
$ws.Range("I89").Value = 843.36365
$ws.Range("J89").Value = 827.6
$ws.Range("K89").Value = 4216.81825
$ws.Range("L89").Value = 4138
$ws.Range("M89").Value = 1399.18175
$ws.Range("N89").Value = -15370

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 28450
$ws.Range("J90").Value = 28450
$ws.Range("L90").Value = 85350
$ws.Range("N90").Value = -97830

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1891.5
$ws.Range("J106").Value = 1700
$ws.Range("L106").Value = 1700
$ws.Range("N106").Value = -2962

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2125
$ws.Range("J116").Value = 2250
$ws.Range("L116").Value = 2250
$ws.Range("N116").Value = -9134

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1079.421
$ws.Range("I132").Value = 1088.7646
$ws.Range("K132").Value = 3266.2938
$ws.Range("M132").Value = -736.2937999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 54942.855
$ws.Range("J134").Value = 54942.855
$ws.Range("L134").Value = 54942.855
$ws.Range("N134").Value = -65082.855

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 77332
$ws.Range("J140").Value = 77332
$ws.Range("L140").Value = 77332
$ws.Range("N140").Value = -87692

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5659.3335
$ws.Range("I63").Value = 5659.3335
$ws.Range("K63").Value = 5659.3335
$ws.Range("M63").Value = -4973.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 5659.3335
$ws.Range("I66").Value = 5659.3335
$ws.Range("K66").Value = 28296.6675
$ws.Range("M66").Value = -24864.6675

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1843.6666
$ws.Range("I132").Value = 1449.25
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 4347.75
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -1817.75
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1027.8334
$ws.Range("I99").Value = 1027.8334
$ws.Range("K99").Value = 1027.8334
$ws.Range("M99").Value = 470.1666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2095.6155
$ws.Range("I105").Value = 2095.1738
$ws.Range("K105").Value = 2095.1738
$ws.Range("M105").Value = -348.1738

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 44833.332
$ws.Range("J68").Value = 44833.332
$ws.Range("L68").Value = 44833.332
$ws.Range("N68").Value = -46331.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 44833.332
$ws.Range("J71").Value = 44833.332
$ws.Range("L71").Value = 134499.996
$ws.Range("N71").Value = -141987.996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1198.25
$ws.Range("I122").Value = 1265.3334
$ws.Range("J122").Value = 997
$ws.Range("K122").Value = 3796.0002
$ws.Range("L122").Value = 2991
$ws.Range("M122").Value = -1346.0002
$ws.Range("N122").Value = -7891

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2182.5789
$ws.Range("I132").Value = 1675.5
$ws.Range("K132").Value = 5026.5
$ws.Range("M132").Value = -2496.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 200001440
$ws.Range("I17").Value = 500000600
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 1500001800
$ws.Range("L17").Value = 6000
$ws.Range("M17").Value = -1500001631
$ws.Range("N17").Value = -6338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2332.6667
$ws.Range("J39").Value = 2332.6667
$ws.Range("L39").Value = 6998.000100000001
$ws.Range("N39").Value = -7586.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 6000
$ws.Range("J104").Value = 6000
$ws.Range("L104").Value = 18000
$ws.Range("N104").Value = -23242

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1470.85
$ws.Range("J107").Value = 1769.6923
$ws.Range("L107").Value = 5309.0769
$ws.Range("N107").Value = -9149.0769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 3436.4285
$ws.Range("I109").Value = 1264.5
$ws.Range("J109").Value = 6332.3335
$ws.Range("K109").Value = 3793.5
$ws.Range("L109").Value = 18997.0005
$ws.Range("M109").Value = -2753.5
$ws.Range("N109").Value = -21077.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 2560
$ws.Range("I112").Value = 1900
$ws.Range("K112").Value = 5700
$ws.Range("M112").Value = -4592

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 17859950
$ws.Range("J114").Value = 28575434
$ws.Range("L114").Value = 85726302
$ws.Range("N114").Value = -85732810

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 41727100
$ws.Range("J131").Value = 144361.8
$ws.Range("L131").Value = 433085.4
$ws.Range("N131").Value = -443165.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 3752434.5
$ws.Range("I7").Value = 5000000
$ws.Range("J7").Value = 1673158.4
$ws.Range("K7").Value = 5000000
$ws.Range("L7").Value = 1673158.4
$ws.Range("M7").Value = -4999888
$ws.Range("N7").Value = -1673382.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H8").Value = 3752434.5
$ws.Range("I8").Value = 5000000
$ws.Range("J8").Value = 1673158.4
$ws.Range("K8").Value = 5000000
$ws.Range("L8").Value = 1673158.4
$ws.Range("M8").Value = -4999861
$ws.Range("N8").Value = -1673436.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1927105.4
$ws.Range("I132").Value = 3208210.2
$ws.Range("J132").Value = 5448
$ws.Range("K132").Value = 9624630.600000001
$ws.Range("L132").Value = 16344
$ws.Range("M132").Value = -9622100.600000001
$ws.Range("N132").Value = -21404

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 49582.445
$ws.Range("J135").Value = 49582.445
$ws.Range("L135").Value = 49582.445
$ws.Range("N135").Value = -59722.445

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 31949.5
$ws.Range("J141").Value = 31949.5
$ws.Range("L141").Value = 31949.5
$ws.Range("N141").Value = -42309.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1378.8462
$ws.Range("I132").Value = 699.75
$ws.Range("J132").Value = 2465.4
$ws.Range("K132").Value = 2099.25
$ws.Range("L132").Value = 7396.200000000001
$ws.Range("M132").Value = 430.75
$ws.Range("N132").Value = -12456.2
